$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A67").Value = "2025/12/05 13:00"
$ws.Range("B67").Value = "-"
$ws.Range("C67").Value = "-"
$ws.Range("D67").Value = "-"
$ws.Range("E67").Value = "-"
$ws.Range("F67").Value = "-"
$ws.Range("G67").Value = "-"
